# Daily scrape update - 2025-10-30
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (C, D, F, H) ---
# ColumnWidth values below are chosen so the persisted OOXML `width`
# (in "Maximum Digit Width" character units) comes out to the exact
# integer targets: C=61, D=70, F=16, H=22.
$ws.Columns.Item(3).ColumnWidth = 60.15625
$ws.Columns.Item(4).ColumnWidth = 69.140625
$ws.Columns.Item(6).ColumnWidth = 15.234375
$ws.Columns.Item(8).ColumnWidth = 21.09375

# --- Column A (OPPORTUNITY ID) holds numeric-looking IDs that must stay text ---
$ws.Range("A2:A6").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "1328856"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328856"
$ws.Range("C2").Value = "Accelerate Romania | Digital Marketing & Community Manager"
$ws.Range("D2").Value = "Iași, Romania"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Nouveaux"

# --- Row 3 ---
$ws.Range("A3").Value = "1328831"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328831"
$ws.Range("C3").Value = "Materials Researcher"
$ws.Range("D3").Value = "Wageningen, Nederland"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "3 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "SHR Timber research"

# --- Row 4 ---
$ws.Range("A4").Value = "1328828"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328828"
$ws.Range("C4").Value = "Laboratory Researcher"
$ws.Range("D4").Value = "Wageningen, Nederland"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "2 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "SHR Timber research"

# --- Row 5 ---
$ws.Range("A5").Value = "1328558"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328558"
$ws.Range("C5").Value = "Flutter Developer"
$ws.Range("D5").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Techno square"

# --- Row 6 ---
$ws.Range("A6").Value = "1321823"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1321823"
$ws.Range("C6").Value = "Sales Responsible at OnurPlas"
$ws.Range("D6").Value = "Konya, Türkiye"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "82 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Onur Plastic"

# --- Remove rows 7-10 (old rows no longer present in the latest scrape) ---
$ws.Range("A7:H10").Delete()
